$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D10:D13").Value = "0.225 (µg/ml)"
$ws.Range("D14:D17").Value = "0.45 (µg/ml)"
$ws.Range("D18:D21").Value = "0.9 (µg/ml)"
